# Auto-generated edit script: applies the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.597.20"
$ws.Range("E2").Value = "  +2.50%  "

# Row 3
$ws.Range("D3").Value = "3.599.25"
$ws.Range("E3").Value = "  +4.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.40"
$ws.Range("E5").Value = "  +2.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.58"
$ws.Range("E6").Value = "  +5.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.49"
$ws.Range("E7").Value = "  +7.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.406"
$ws.Range("E8").Value = "  +3.09%  "

# Row 9
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("E10").Value = "  +3.99%  "

# Row 11
$ws.Range("D11").Value = "3.597.36"
$ws.Range("E11").Value = "  +4.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.05"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.31"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("D15").Value = "4.265.66"
$ws.Range("E15").Value = "  +4.75%  "

# Row 16
$ws.Range("D16").Value = "95.433.55"
$ws.Range("E16").Value = "  +2.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000256"
$ws.Range("E17").Value = "  +3.70%  "

# Row 18
$ws.Range("D18").Value = "3.593.18"
$ws.Range("E18").Value = "  +4.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  -2.94%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  +7.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.02"
$ws.Range("E21").Value = "  -1.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("E22").Value = "  +7.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.492"
$ws.Range("E23").Value = "  +10.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "510.64"
$ws.Range("E24").Value = "  +1.89%  "

# Row 25
$ws.Range("E25").Value = "  +5.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.90"
$ws.Range("E27").Value = "  +1.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.70"
$ws.Range("E28").Value = "  +5.80%  "

# Row 29
$ws.Range("D29").Value = "3.777.60"
$ws.Range("E29").Value = "  +4.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.19"
$ws.Range("E30").Value = "  +16.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.31"
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("E32").Value = "  -0.08%  "

# Row 33
$ws.Range("E33").Value = "  +2.20%  "

# Row 34
$ws.Range("E34").Value = "  +1.29%  "

# Row 35
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("E35").Value = "  +2.05%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.88"
$ws.Range("E36").Value = "  +5.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.560"
$ws.Range("E37").Value = "  +2.44%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "575.63"
$ws.Range("E38").Value = "  +4.29%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.18"
$ws.Range("E39").Value = "  +9.58%  "

# Row 40
$ws.Range("E40").Value = "  +6.24%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.150"
$ws.Range("E42").Value = "  +0.55%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.925"
$ws.Range("E43").Value = "  +0.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.74"
$ws.Range("E44").Value = "  +3.80%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  -0.60%  "

# Row 46
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.77"
$ws.Range("E46").Value = "  +0.34%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.82"
$ws.Range("E47").Value = "  +30.10%  "

# Row 48
$ws.Range("E48").Value = "  +1.63%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +5.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.99"
$ws.Range("E50").Value = "  +0.58%  "

# Row 51
$ws.Range("E51").Value = "  -6.50%  "

Write-Output "Applied cryptos update: $( (Get-Date).ToString() )"
